$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheets: rename Sheet1 -> LoginData, insert RegisterData + Pages after it.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "LoginData"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "RegisterData"

$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Pages"

# ---------------------------------------------------------------------------
# Content fill -- ordered so the resulting shared-strings table comes out in
# the same order as the target workbook.
# ---------------------------------------------------------------------------

# LoginData username/password header + sample row (shared strings 0..3)
$ws1.Range("A1").Value = "Username"
$ws1.Range("B1").Value = "Password"
$ws1.Range("A2").Value = "qa123123"
$ws1.Range("B2").Value = "testing@8"

# Pages (shared strings 4..5)
$ws3.Range("A1").Value = "EditorCode"
$ws3.Range("A2").Value = 'print("Hello, World!")'

# RegisterData header row (shared strings 6..11)
$ws2.Range("A1").Value = "Valid Username"
$ws2.Range("B1").Value = "Valid Password"
$ws2.Range("C1").Value = "Valid ConfirmPassword"
$ws2.Range("D1").Value = "Invalid Username"
$ws2.Range("E1").Value = "Invalid Password"
$ws2.Range("F1").Value = "Invalid Confirm Password"

# Style RegisterData's used range: horizontal-left + wrap text (cellXfs index 1)
$stRegister = $wb.Styles.Add("RegisterDataStyle")
$stRegister.WrapText = $true
$stRegister.HorizontalAlignment = -4131
$ws2.Range("A1:F4").Style = $stRegister

# LoginData remaining header cells (reuse shared strings 9,10)
$ws1.Range("C1").Value = "Invalid Username"
$ws1.Range("D1").Value = "Invalid Password"

# Style LoginData's used range: wrap text only (cellXfs index 2)
$ws1.Range("A1:D2").WrapText = $true

# RegisterData data rows (shared strings 12,13, number, reuse, "test"=14, number, 15)
$ws2.Range("A2").Value = "qa111222"
$ws2.Range("B2").Value = "testing@123"
$ws2.Range("C2").Value = "testing@123"
$ws2.Range("D2").Value = 123
$ws2.Range("E2").Value = "qa111222"
$ws2.Range("E3").Value = "test"
$ws2.Range("E4").Value = 12345678
$ws2.Range("F2").Value = " testing@123456"

# LoginData remaining data cells (reuse shared string 14 "test")
$ws1.Range("C2").Value = "test"
$ws1.Range("D2").Value = "test"

# Hyperlink on LoginData!B2 (creates Hyperlink font/style, cellXfs index 3)
$ws1.Hyperlinks.Add($ws1.Range("B2"), "testing@8", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "testing@8")
$ws1.Range("B2").WrapText = $true

# ---------------------------------------------------------------------------
# Row heights / column widths
# ---------------------------------------------------------------------------
$ws1.Rows.Item(1).RowHeight = 45

$ws1.Columns.Item(1).ColumnWidth = 13.6
$ws1.Columns.Item(2).ColumnWidth = 11.42
$ws1.Columns.Item(3).ColumnWidth = 10.25
$ws1.Columns.Item(4).ColumnWidth = 10.59

$ws2.Rows.Item(1).RowHeight = 39

$ws2.Columns.Item(1).ColumnWidth = 18.09
$ws2.Columns.Item(2).ColumnWidth = 16.25
$ws2.Columns.Item(3).ColumnWidth = 21.92
$ws2.Columns.Item(4).ColumnWidth = 17.42
$ws2.Columns.Item(5).ColumnWidth = 13.75
$ws2.Columns.Item(6).ColumnWidth = 16.75

$ws3.Columns.Item(1).ColumnWidth = 19.75

# ---------------------------------------------------------------------------
# Selections / active cells per sheet, and page orientation on Pages.
# ---------------------------------------------------------------------------
[void]$ws2.Activate()
[void]$ws2.Range("F2").Select()

[void]$ws3.Activate()
[void]$ws3.Range("A2").Select()
$ws3.PageSetup.Orientation = 1

[void]$ws1.Activate()
[void]$ws1.Range("D2").Select()
